# "sua them excel image vao chi tiet san pham"
# - Replace the six embedded pictures on the "chi tiet sp" sheet with plain
#   text file paths typed into the "anh mau sac" / "anh chinh" columns.
# - Drop the now-redundant second data row (id 25) from that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("chi tiet sp")
$ws.Activate() | Out-Null

# Remove every picture shape that used to be drawn over this sheet.
while ($ws.Shapes.Count -gt 0) {
    $ws.Shapes.Item(1).Delete() | Out-Null
}

# Row 3 (the second product-detail record) is removed entirely.
$ws.Rows.Item(3).Delete() | Out-Null

# Column M ("anh chinh") gets a single image path; column L ("anh mau sac")
# gets a comma-separated list of image paths. Set M2 first so the new
# shared-string table keeps the single path ahead of the multi-path one.
$ws.Range("M2").Value = "D:\anh\digital-art-artwork-illustration-digital-painting-pink-hd-wallpaper-cc923972d4ec819cccf9880f8b916d69.jpg"
$ws.Range("L2").Value = "D:\anh\318500342_195400319734317_2076385213486251926_n.jpg,D:\anh\28.png,D:\anh\18.jpg,D:\anh\318500342_195400319734317_2076385213486251926_n.jpg"

# Restore the on-screen selection to what the author left it at.
$ws.Range("K12").Select() | Out-Null
